$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 117, shifting existing rows 117..197 down to 118..198
$ws.Rows.Item(117).Insert()

# Populate the newly inserted row 117 with the new record
$ws.Cells.Item(117, 1).Value = 10
$ws.Cells.Item(117, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(117, 3).Value = "La Araucanía"
$ws.Cells.Item(117, 4).Value = 45167
$ws.Cells.Item(117, 4).NumberFormat = $ws.Cells.Item(118, 4).NumberFormat
$ws.Cells.Item(117, 5).Value = 9
$ws.Cells.Item(117, 6).Value = "Fruta"
$ws.Cells.Item(117, 7).Value = 100107
$ws.Cells.Item(117, 8).Value = "Otros"
$ws.Cells.Item(117, 9).Value = 100107002
$ws.Cells.Item(117, 10).Value = "Chirimoya"
$ws.Cells.Item(117, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(117, 12).Value = "Primera"
$ws.Cells.Item(117, 13).Value = 35
$ws.Cells.Item(117, 14).Value = 3300
$ws.Cells.Item(117, 15).Value = 3300
$ws.Cells.Item(117, 16).Value = 3300
$ws.Cells.Item(117, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(117, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(117, 19).Value = 3300
$ws.Cells.Item(117, 20).Value = 1
